$wb = $excel.ActiveWorkbook

# --- info_solution sheet: update comp_time (A2) ---
$wsInfo = $wb.Worksheets.Item("info_solution")
$wsInfo.Range("A2").Value = 1.3165929317474365

# --- design_users sheet: update Peak demand [kW] (B) and Yearly Demand [MWh] (C) ---
$wsDesign = $wb.Worksheets.Item("design_users")

# user1 (row 2)
$wsDesign.Range("B2").Value = 30523.270014740006
$wsDesign.Range("C2").Value = 232.10403242918812

# user2 (row 3)
$wsDesign.Range("B3").Value = 12160.103553284996
$wsDesign.Range("C3").Value = 92.467454113238

# user3 (row 4)
$wsDesign.Range("B4").Value = 24500.351344201987
$wsDesign.Range("C4").Value = 186.30475503361959
